$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(37, "erwrewrw", "04-11-2023"),
    @(38, "erwerwrew", "04-11-2023"),
    @(39, "dfdsfsfs", "04-11-2023"),
    @(40, "erwrwerw", "04-11-2023")
)

$startRow = 39
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    # Force column C to text so date-looking strings ("04-11-2023") are
    # stored verbatim (t="str"/shared-string) instead of being parsed
    # into a serial date number, matching the existing rows above.
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $data[2]
}

# The workbook originally carries an "ignore number-stored-as-text"
# marker over A1:C38; extend it to cover the newly appended rows
# (A1:C42) the same way Excel does when you dismiss/ignore the warning
# on the enlarged range.
$ws.Range("A1:C42").Errors.Item(3).Ignore = $true
